$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 51
$ws.Range("C3").Value = 55
$ws.Range("C4").Value = 52
$ws.Range("C5").Value = 54
$ws.Range("C6").Value = 54
$ws.Range("C7").Value = 54
$ws.Range("C8").Value = 53
$ws.Range("C9").Value = 56
$ws.Range("C10").Value = 51
$ws.Range("C11").Value = 49
$ws.Range("C12").Value = 43
$ws.Range("C13").Value = 50
$ws.Range("C14").Value = 52

$ws.Range("C15").Select()
